$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Input" (sheet1): fill in the digit sequence across row 1 (C1:Q1)
# and the start of row 2 (A2:C2), clear out the placeholder "1"s that were
# typed across row 6, and move the active selection to C3.
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

$row1Cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$row1Vals = @(1,2,3,4,5,6,7,8,9,0,1,2,3,4,5)
for ($i = 0; $i -lt $row1Cols.Length; $i++) {
    $wsInput.Range($row1Cols[$i] + "1").Value = $row1Vals[$i]
}

$wsInput.Range("A2").Value = 6
$wsInput.Range("B2").Value = 7
$wsInput.Range("C2").Value = 8

$wsInput.Range("A6:Q6").ClearContents()

$wsInput.Range("C3").Select()

# ---------------------------------------------------------------------------
# Sheet "View_Print" (sheet2): give A1:B1 an explicit (general) horizontal
# alignment, merge them into a single cell, and move the selection to D2.
# ---------------------------------------------------------------------------
$wsPrint = $wb.Worksheets.Item("View_Print")

$wsPrint.Range("A1:B1").HorizontalAlignment = 1
$wsPrint.Range("A1:B1").Merge()

$wsPrint.Range("D2").Select()

# Restore "Input" as the active/selected tab (it was the active sheet before
# these edits and stays that way afterwards).
$wsInput.Activate()

Write-Host "edit applied"
